$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item(1)
# Row 42
$ws.Range("H42").Value = 56.142857
$ws.Range("I42").Value = 42.8
$ws.Range("J42").Value = 89.5
$ws.Range("K42").Value = 128.4
$ws.Range("L42").Value = 268.5
$ws.Range("M42").Value = 101.6
$ws.Range("N42").Value = -728.5

# Row 76
$ws.Range("H76").Value = 71434060
$ws.Range("I76").Value = 4364.2
$ws.Range("K76").Value = 4364.2
$ws.Range("M76").Value = -4049.2

# Row 79
$ws.Range("H79").Value = 71434060
$ws.Range("I79").Value = 4364.2
$ws.Range("K79").Value = 4364.2
$ws.Range("M79").Value = -3272.2

# Row 132
$ws.Range("H132").Value = 3131.8096
$ws.Range("I132").Value = 3383.5789
$ws.Range("K132").Value = 10150.7367
$ws.Range("M132").Value = -7620.736699999999

# ---- ARM ----
$ws = $wb.Worksheets.Item(2)
# Row 32
$ws.Range("H32").Value = 3140.0144
$ws.Range("I32").Value = 3184.0293
$ws.Range("K32").Value = 3184.0293
$ws.Range("M32").Value = -2897.0293

# Row 132
$ws.Range("H132").Value = 3907.0908
$ws.Range("I132").Value = 3997.9524
$ws.Range("K132").Value = 11993.8572
$ws.Range("M132").Value = -9463.8572

# ---- BSM ----
$ws = $wb.Worksheets.Item(3)
# Row 20
$ws.Range("H20").Value = 3068.5
$ws.Range("I20").Value = 2626.7058
$ws.Range("K20").Value = 2626.7058
$ws.Range("M20").Value = -2379.7058

# Row 37
$ws.Range("H37").Value = 2311.75
$ws.Range("I37").Value = 1645.2
$ws.Range("J37").Value = 2533.9333
$ws.Range("K37").Value = 1645.2
$ws.Range("L37").Value = 2533.9333
$ws.Range("M37").Value = -1508.2
$ws.Range("N37").Value = -2807.9333

# Row 82
$ws.Range("H82").Value = 53600
$ws.Range("I82").Value = 9333.333000000001
$ws.Range("K82").Value = 9333.333000000001
$ws.Range("M82").Value = -8950.333000000001

# Row 85
$ws.Range("H85").Value = 53600
$ws.Range("I85").Value = 9333.333000000001
$ws.Range("K85").Value = 9333.333000000001
$ws.Range("M85").Value = -8007.333000000001

# Row 86
$ws.Range("H86").Value = 851662.8
$ws.Range("I86").Value = 1001579.7
$ws.Range("J86").Value = 2133.6667
$ws.Range("K86").Value = 1001579.7
$ws.Range("L86").Value = 2133.6667
$ws.Range("M86").Value = -1000456.7
$ws.Range("N86").Value = -4379.6667

# Row 89
$ws.Range("H89").Value = 851662.8
$ws.Range("I89").Value = 1001579.7
$ws.Range("J89").Value = 2133.6667
$ws.Range("K89").Value = 5007898.5
$ws.Range("L89").Value = 10668.3335
$ws.Range("M89").Value = -5002282.5
$ws.Range("N89").Value = -21900.3335

# Row 97
$ws.Range("H97").Value = 14899.8
$ws.Range("I97").Value = 6125
$ws.Range("K97").Value = 6125
$ws.Range("M97").Value = -5134

# ---- CRP ----
$ws = $wb.Worksheets.Item(4)
# Row 14
$ws.Range("H14").Value = 6121.4287
$ws.Range("I14").Value = 1500
$ws.Range("J14").Value = 6891.6665
$ws.Range("K14").Value = 1500
$ws.Range("L14").Value = 6891.6665
$ws.Range("M14").Value = -1330
$ws.Range("N14").Value = -7231.6665

# Row 62
$ws.Range("H62").Value = 6500
$ws.Range("I62").Value = 5000
$ws.Range("J62").Value = 8000
$ws.Range("K62").Value = 5000
$ws.Range("L62").Value = 8000
$ws.Range("M62").Value = -4376
$ws.Range("N62").Value = -9248

# Row 65
$ws.Range("H65").Value = 6500
$ws.Range("I65").Value = 5000
$ws.Range("J65").Value = 8000
$ws.Range("K65").Value = 25000
$ws.Range("L65").Value = 40000
$ws.Range("M65").Value = -21880
$ws.Range("N65").Value = -46240

# Row 132
$ws.Range("H132").Value = 1466
$ws.Range("I132").Value = 1288
$ws.Range("K132").Value = 3864
$ws.Range("M132").Value = -1334

# Row 141
$ws.Range("H141").Value = 493012.84
$ws.Range("J141").Value = 555132.3
$ws.Range("L141").Value = 555132.3
$ws.Range("N141").Value = -565492.3

# ---- CUL ----
$ws = $wb.Worksheets.Item(5)
# Row 2
$ws.Range("H2").Value = 243.46153
$ws.Range("I2").Value = 54
$ws.Range("K2").Value = 324
$ws.Range("M2").Value = -211

# Row 6
$ws.Range("H6").Value = 120.47369
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()

# Row 37
$ws.Range("H37").Value = 117756.43
$ws.Range("J37").Value = 117756.43
$ws.Range("L37").Value = 353269.29
$ws.Range("N37").Value = -353493.29

# ---- GSM ----
$ws = $wb.Worksheets.Item(6)
# Row 17
$ws.Range("H17").Value = 3617.4167
$ws.Range("J17").Value = 3936.4546
$ws.Range("L17").Value = 3936.4546
$ws.Range("N17").Value = -4272.4546

# Row 26
$ws.Range("H26").Value = 24219.5
$ws.Range("J26").Value = 24219.5
$ws.Range("L26").Value = 24219.5
$ws.Range("N26").Value = -24779.5

# Row 43
$ws.Range("H43").Value = 20003
$ws.Range("I43").Value = 10004
$ws.Range("K43").Value = 10004
$ws.Range("M43").Value = -9853

# Row 50
$ws.Range("H50").Value = 24219.5
$ws.Range("J50").Value = 24219.5
$ws.Range("L50").Value = 24219.5
$ws.Range("N50").Value = -25215.5

# Row 52
$ws.Range("H52").Value = 28499.666
$ws.Range("J52").Value = 28499.666
$ws.Range("L52").Value = 28499.666
$ws.Range("N52").Value = -29017.666

# Row 70
$ws.Range("H70").Value = 8683.857
$ws.Range("I70").Value = 7944.55
$ws.Range("J70").Value = 10532.125
$ws.Range("K70").Value = 7944.55
$ws.Range("L70").Value = 10532.125
$ws.Range("M70").Value = -7674.55
$ws.Range("N70").Value = -11072.125

# Row 73
$ws.Range("H73").Value = 8683.857
$ws.Range("I73").Value = 7944.55
$ws.Range("J73").Value = 10532.125
$ws.Range("K73").Value = 7944.55
$ws.Range("L73").Value = 10532.125
$ws.Range("M73").Value = -7008.55
$ws.Range("N73").Value = -12404.125

# Row 97
$ws.Range("H97").Value = 755
$ws.Range("I97").Value = 1010
$ws.Range("J97").Value = 500
$ws.Range("K97").Value = 1010
$ws.Range("L97").Value = 500
$ws.Range("M97").Value = -514
$ws.Range("N97").Value = -1492

# Row 102
$ws.Range("H102").Value = 837.4737
$ws.Range("I102").Value = 609.625
$ws.Range("K102").Value = 609.625
$ws.Range("M102").Value = 1012.375

# Row 122
$ws.Range("H122").Value = 4074.9167
$ws.Range("I122").Value = 1899.8334
$ws.Range("J122").Value = 6250
$ws.Range("K122").Value = 5699.5002
$ws.Range("L122").Value = 18750
$ws.Range("M122").Value = -3249.5002
$ws.Range("N122").Value = -23650

# ---- LTW ----
$ws = $wb.Worksheets.Item(7)
# Row 46
$ws.Range("H46").Value = 3906.6333
$ws.Range("I46").Value = 3387.96
$ws.Range("J46").Value = 6500
$ws.Range("K46").Value = 3387.96
$ws.Range("L46").Value = 6500
$ws.Range("M46").Value = -3199.96
$ws.Range("N46").Value = -6876

# Row 55
$ws.Range("H55").Value = 1275.9375
$ws.Range("I55").Value = 334.1111
$ws.Range("J55").Value = 2486.8572
$ws.Range("K55").Value = 334.1111
$ws.Range("L55").Value = 2486.8572
$ws.Range("M55").Value = -161.1111
$ws.Range("N55").Value = -2832.8572

# Row 122
$ws.Range("H122").Value = 2225768.2
$ws.Range("I122").Value = 1669335
$ws.Range("J122").Value = 3338635
$ws.Range("K122").Value = 5008005
$ws.Range("L122").Value = 10015905
$ws.Range("M122").Value = -5005555
$ws.Range("N122").Value = -10020805

# ---- WVR ----
$ws = $wb.Worksheets.Item(8)
# Row 21
$ws.Range("H21").Value = 0
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()

# Row 35
$ws.Range("H35").Value = 0
$ws.Range("J35").Value = 0
$ws.Range("L35").Value = 0
$ws.Range("N35").ClearContents()

# Row 94
$ws.Range("H94").Value = 46247.5
$ws.Range("I94").Value = 40596.2
$ws.Range("J94").Value = 55666.332
$ws.Range("K94").Value = 40596.2
$ws.Range("L94").Value = 55666.332
$ws.Range("M94").Value = -39695.2
$ws.Range("N94").Value = -57468.332

# Row 122
$ws.Range("H122").Value = 50002144
$ws.Range("I122").Value = 62501476
$ws.Range("J122").Value = 4826.25
$ws.Range("K122").Value = 187504428
$ws.Range("L122").Value = 14478.75
$ws.Range("M122").Value = -187501978
$ws.Range("N122").Value = -19378.75

# Row 123
$ws.Range("H123").Value = 64999.5
$ws.Range("J123").Value = 64999.5
$ws.Range("L123").Value = 64999.5
$ws.Range("N123").Value = -74799.5
